$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" before the current "2022-Q2" sheet
#    (which sits right after "总计" at position 2).
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q3"

# Header row
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# Data rows.
# NOTE: columns B (fund code) and D-G (numeric-looking figures kept as
# *text* in the source data, per the diff's t="inlineStr") would otherwise
# be auto-coerced to numbers by Excel's smart-typing (e.g. "005396" ->
# 5396, losing the leading zero). Prefixing with an apostrophe forces
# Excel to keep them as literal text, exactly like a user typing
# '005396 into a General cell; re-applying the "Normal" style afterwards
# drops the resulting quote-prefix flag so the cell ends up plain text
# with no special formatting (matching the source sheet).
$newSheet.Range("A2").Value2 = 0
$newSheet.Range("B2").Value2 = "'005396"
$newSheet.Range("C2").Value2 = "中金丰硕混合"
$newSheet.Range("D2").Value2 = "'1.77"
$newSheet.Range("E2").Value2 = "'76.61"
$newSheet.Range("F2").Value2 = "'7.64"
$newSheet.Range("G2").Value2 = "'0.1352"
$newSheet.Range("H2").Value2 = 2

$newSheet.Range("A3").Value2 = 1
$newSheet.Range("B3").Value2 = "'013659"
$newSheet.Range("C3").Value2 = "中融金融鑫选3个月持有混合A"
$newSheet.Range("D3").Value2 = "'1.32"
$newSheet.Range("E3").Value2 = "'85.15"
$newSheet.Range("F3").Value2 = "'5.01"
$newSheet.Range("G3").Value2 = "'0.0661"
$newSheet.Range("H3").Value2 = 8

$newSheet.Range("A4").Value2 = 2
$newSheet.Range("B4").Value2 = "'013660"
$newSheet.Range("C4").Value2 = "中融金融鑫选3个月持有混合C"
$newSheet.Range("D4").Value2 = "'0.81"
$newSheet.Range("E4").Value2 = "'85.15"
$newSheet.Range("F4").Value2 = "'5.01"
$newSheet.Range("G4").Value2 = "'0.0406"
$newSheet.Range("H4").Value2 = 8

$newSheet.Range("A5").Value2 = 3
$newSheet.Range("B5").Value2 = "'516980"
$newSheet.Range("C5").Value2 = "华富中证证券公司先锋策略ETF"
$newSheet.Range("D5").Value2 = "'0.28"
$newSheet.Range("E5").Value2 = "'99.02"
$newSheet.Range("F5").Value2 = "'3.04"
$newSheet.Range("G5").Value2 = "'0.0085"
$newSheet.Range("H5").Value2 = 7

# Drop the quote-prefix flag picked up from the apostrophe-forced text
# entry above so the cells end up unstyled, like the source data.
$newSheet.Range("B2:B5").Style = "Normal"
$newSheet.Range("D2:G5").Style = "Normal"

# Header cells (B1:H1) use the bold "header" style -- copy it from the
# equivalent header on the following sheet (2022-Q2), which already has
# the right formatting, then re-apply the data-row values above stay intact
# because PasteSpecial(formats) does not touch cell contents.
$otherHeader = $wb.Worksheets.Item("2022-Q2").Range("B1:H1")
$otherHeader.Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$otherIndexCol = $wb.Worksheets.Item("2022-Q2").Range("A2")
$otherIndexCol.Copy() | Out-Null
$newSheet.Range("A2:A5").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: push the existing rows down by
#    one and insert the new 2022-Q3 row at the top of the data block.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 8; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $summary.Range("B$dst").Value2 = $summary.Range("B$src").Value2
    $summary.Range("C$dst").Value2 = $summary.Range("C$src").Value2
    $summary.Range("D$dst").Value2 = $summary.Range("D$src").Value2
}

# Row 9 is brand new -- copy the style of row 8's index cell (A) onto it.
$summary.Range("A8").Copy() | Out-Null
$summary.Range("A9").PasteSpecial(-4122) | Out-Null

# Column A is a simple 0-based row index; recompute it for every data row.
for ($r = 2; $r -le 9; $r++) {
    $summary.Range("A$r").Value2 = $r - 2
}

# Row 2 becomes the new 2022-Q3 entry.
$summary.Range("B2").Value2 = "2022-Q3"
$summary.Range("C2").Value2 = 4
$summary.Range("D2").Value2 = 0.25

Write-Host "edit complete"
